$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("L2").Value = 0.91

# Row 3
$ws1.Range("L3").Value = 1.11

# Row 4
$ws1.Range("D4").Value = 1
$ws1.Range("H4").Value = 0.79
$ws1.Range("I4").Value = "Low"
$ws1.Range("L4").Value = 0.9

# Row 5
$ws1.Range("D5").Value = 1
$ws1.Range("L5").Value = 1.2

# Row 6
$ws1.Range("D6").Value = 1
$ws1.Range("L6").Value = 1.04

# Row 7
$ws1.Range("D7").Value = 1
$ws1.Range("L7").Value = 1.09

# Row 8
$ws1.Range("D8").Value = 1
$ws1.Range("L8").Value = 1

# Row 9
$ws1.Range("D9").Value = 1
$ws1.Range("L9").Value = 0.9

# Row 10
$ws1.Range("D10").Value = 1
$ws1.Range("L10").Value = 1.13

# Row 11
$ws1.Range("D11").Value = 2
$ws1.Range("L11").Value = 1.19

# Row 12
$ws1.Range("D12").Value = 1
$ws1.Range("L12").Value = 1.12

# Row 13
$ws1.Range("D13").Value = 1
$ws1.Range("L13").Value = 0.84

# Row 14
$ws1.Range("D14").Value = 1
$ws1.Range("L14").Value = 0.9

# Row 15
$ws1.Range("D15").Value = 1
$ws1.Range("L15").Value = 0.88

# Row 16
$ws1.Range("D16").Value = 1
$ws1.Range("L16").Value = 1.01

# Row 17
$ws1.Range("D17").Value = 1
$ws1.Range("L17").Value = 1.15

# --- Sheet: Summary ---
# These cells hold numeric-looking values but are stored as text in the
# workbook, so a leading apostrophe is used to force text entry (same as
# typing '31 directly into Excel) rather than letting Excel infer a number.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").Value = "'31"
$ws2.Range("B10").Value = "'16"
$ws2.Range("B11").Value = "'8"
$ws2.Range("B12").Value = "'3"
